$d = $word.ActiveDocument

# --- 1. Add the "Victim's Attorney" line to the "Copies served ... to:" ---
#        paragraph, right after the existing "County Jail: PS   EM;" text.

$apos = [char]0x2019
$newLine = "Victim" + $apos + "s Attorney (if applicable): PS   OS   EM"

$rng = $d.Content
$found = $rng.Find.Execute("County Jail: PS   EM;", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'County Jail: PS   EM;' run to anchor the new text on."
}

# Collapse to the end of the match (right before the paragraph mark) and add
# the new content as two runs: a separating space, then the new label line.
$rng.Collapse(0)
$insertStart = $rng.Start

$rng.InsertAfter(" ")
$spaceRun = $d.Range($insertStart, $insertStart + 1)
$spaceRun.Font.Name = "Palatino Linotype"
$spaceRun.Font.Size = 8

$textStart = $insertStart + 1
$textInsertPoint = $d.Range($textStart, $textStart)
$textInsertPoint.InsertAfter($newLine)
$textRun = $d.Range($textStart, $textStart + $newLine.Length)
$textRun.Font.Name = "Palatino Linotype"
$textRun.Font.Size = 8
